# Updated MCH102 to MCH251 -- add the MCH214 collection record as row 2
# of the finding-aid sheet (header row 1 is left untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data row -----------------------------------------------------
# Columns: A identifier | B alternativeIdentifiers | C title | D date_s |
#          E levelOfDescription | F extentAndMedium | G notes | H file_path
$ws.Range("A2").Value = "MCH214"
$ws.Range("C2").Value = "HOFFNUNGSWANDERUNG FUR SUDAFRIKA"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: CABINET 1C | GRAP COUNT NUMER: NONE"
# B2, D2 and H2 stay empty (alternativeIdentifiers, date_s, file_path
# were not recorded for this item) but still pick up the row's formatting.

# --- Formatting for the new row ---------------------------------------
# Match the sheet's data-row style: 10pt Calibri, automatic/theme text color.
foreach ($addr in @("A2","C2","D2","E2","F2","G2","H2")) {
    $cell = $ws.Range($addr)
    $cell.Font.Name = "Calibri"
    $cell.Font.Size = 10
    $cell.Font.ThemeColor = 1
}

# --- View state: keep header frozen, select the new row ---------------
$ws.Range("A2:H2").Select()
$excel.ActiveWindow.FreezePanes = $true
